$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = [double]"599.8753063347568"
    "C2" = [double]"1199.7506126695166"
    "D2" = [double]"-1.2862974574368025e-12"
    "E2" = [double]"1199.9353063344545"
    "B3" = [double]"599.9678710563109"
    "C3" = [double]"4502071.385664888"
    "D3" = [double]"4500871.449922772"
    "E3" = [double]"1800.0878710555369"
    "B4" = [double]"600.5594735694418"
    "C4" = [double]"8005255.006019249"
    "D4" = [double]"8004053.887072104"
    "E4" = [double]"2400.7394735681664"
    "B5" = [double]"600.2147822294103"
    "C5" = [double]"11253755.7173883"
    "D5" = [double]"11252555.287823813"
    "E5" = [double]"3000.4547822276268"
    "B6" = [double]"600.3214767115759"
    "C6" = [double]"14402697.071327845"
    "D6" = [double]"14401496.428374384"
    "E6" = [double]"3600.6214767092574"
    "B7" = [double]"600.2360796237413"
    "C7" = [double]"17504966.6031143"
    "D7" = [double]"17503766.13095503"
    "E7" = [double]"4200.596079621566"
    "B8" = [double]"599.5091579742073"
    "C8" = [double]"20573210.481114157"
    "D8" = [double]"20572011.4627982"
    "E8" = [double]"4799.929157975466"
    "B9" = [double]"599.8249842873151"
    "C9" = [double]"23629565.619365904"
    "D9" = [double]"23628365.969397295"
    "E9" = [double]"5400.3049842921855"
    "B10" = [double]"600.7954422449061"
    "C10" = [double]"26674163.092877"
    "D10" = [double]"26672961.5019926"
    "E10" = [double]"6001.335442253425"
    "B11" = [double]"600.5838988670773"
    "C11" = [double]"29708261.93013732"
    "D11" = [double]"29707060.762339693"
    "E11" = [double]"6601.183898879313"
    "B12" = [double]"598.549847177036"
    "C12" = [double]"32725534.90517791"
    "D12" = [double]"32724337.80548353"
    "E12" = [double]"7199.209847193052"
    "B13" = [double]"599.2015926794934"
    "C13" = [double]"35747711.194993466"
    "D13" = [double]"35746512.79180807"
    "E13" = [double]"7799.921592699312"
    "B14" = [double]"601.0588600676085"
    "C14" = [double]"38779119.0586033"
    "D14" = [double]"38777916.94088308"
    "E14" = [double]"8401.838860088337"
    "B15" = [double]"600.3978231610255"
    "C15" = [double]"41791607.27100982"
    "D15" = [double]"41790406.47536355"
    "E15" = [double]"9001.237823177104"
    "B16" = [double]"600.7257533291081"
    "C16" = [double]"44805559.540888086"
    "D16" = [double]"44804358.08938161"
    "E16" = [double]"9601.625753340439"
    "B17" = [double]"599.9898066040677"
    "C17" = [double]"47818094.77343345"
    "D17" = [double]"47816894.79382014"
    "E17" = [double]"10200.94980661068"
    "B18" = [double]"600.947707503187"
    "C18" = [double]"50832537.625580296"
    "D18" = [double]"50831335.73016523"
    "E18" = [double]"10801.967707504975"
    "B19" = [double]"600.2994969805065"
    "C19" = [double]"53841989.96946201"
    "D19" = [double]"53840789.37046812"
    "E19" = [double]"11401.379496977459"
    "B20" = [double]"600.3554481919007"
    "C20" = [double]"56851841.82696036"
    "D20" = [double]"56850641.11606393"
    "E20" = [double]"12001.495448184049"
    "B21" = [double]"599.5321962735959"
    "C21" = [double]"59852783.81603763"
    "D21" = [double]"59851584.7516452"
    "E21" = [double]"12600.732196260875"
    "B22" = [double]"600.5676358301349"
    "C22" = [double]"62869930.14398665"
    "D22" = [double]"62868729.00871471"
    "E22" = [double]"13201.827635812604"
    "B23" = [double]"600.8865592689846"
    "C23" = [double]"65876459.28992207"
    "D23" = [double]"65875257.51680335"
    "E23" = [double]"13802.20655924644"
    "B24" = [double]"598.3720865252427"
    "C24" = [double]"68871174.98992705"
    "D24" = [double]"68869978.2457542"
    "E24" = [double]"14399.752086497885"
    "B25" = [double]"600.3936393425197"
    "C25" = [double]"71883989.0505184"
    "D25" = [double]"71882788.26323967"
    "E25" = [double]"15001.833639310173"
    "B26" = [double]"600.7450952211633"
    "C26" = [double]"74889201.71028914"
    "D26" = [double]"74888000.2200993"
    "E26" = [double]"15602.2450951839"
    "B27" = [double]"599.0875700502468"
    "C27" = [double]"77887060.59693499"
    "D27" = [double]"77885862.4217948"
    "E27" = [double]"16200.647570007946"
    "B28" = [double]"600.3245311557905"
    "C28" = [double]"80900262.82426457"
    "D28" = [double]"80899062.17520209"
    "E28" = [double]"16801.944531120753"
    "B29" = [double]"598.4569552037689"
    "C29" = [double]"83893655.02980562"
    "D29" = [double]"83892458.11589517"
    "E29" = [double]"17400.13695518133"
    "B30" = [double]"599.9624324929807"
    "C30" = [double]"86907807.93840191"
    "D30" = [double]"86906608.01353706"
    "E30" = [double]"18001.702432483216"
    "B31" = [double]"599.7492341130744"
    "C31" = [double]"89913343.90883328"
    "D31" = [double]"89912144.41036497"
    "E31" = [double]"18601.549234115864"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
